# Insert a new data row at row 131 of the "Apio" (Vega Modelo de Temuco) sheet.
# This pushes the previously-existing rows 131..260 down to 132..261 (the
# worksheet's used range grows from A1:R260 to A1:R261) and populates the
# newly inserted row 131 with a new price-report record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 131:260 down by one row, inserting a blank row at 131.
$ws.Rows("131:131").Insert()

# Fill in the new row with the new record's data.
$ws.Range("A131").Value = 10
$ws.Range("B131").Value = "Vega Modelo de Temuco"
$ws.Range("C131").Value = "La Araucanía"
$ws.Range("D131").Value = 44629
$ws.Range("E131").Value = 9
$ws.Range("F131").Value = 100112017
$ws.Range("G131").Value = "Apio"
$ws.Range("H131").Value = "Americana (o)"
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 65
$ws.Range("K131").Value = 8000
$ws.Range("L131").Value = 8000
$ws.Range("M131").Value = 8000
$ws.Range("N131").Value = "$/docena de matas"
$ws.Range("O131").Value = "Provincia del Elquí"
$ws.Range("P131").Value = 1333
$ws.Range("Q131").Value = 6
$ws.Range("R131").Value = "Hortaliza"
